# Add a "FiscalYear" parameter column to the Global parameters sheet so the
# fiscal year used for cost calculations can be configured.

$wb = $excel.ActiveWorkbook
$wsGlobal = $wb.Worksheets.Item("Global")
$wsActive = $wb.ActiveSheet

# New header + value in column D.
$wsGlobal.Range("D1").Value = "FiscalYear"
$wsGlobal.Range("D2").Value = 2020

# D2 becomes the new right-hand edge of the bordered parameter block, so it
# should pick up the same box-border formatting that the old last column
# (C2) had; C2 in turn becomes an interior column and loses its right edge.
$wsGlobal.Range("C2").Copy()
$wsGlobal.Range("D2").PasteSpecial(-4122) | Out-Null
$wsGlobal.Range("B2").Copy()
$wsGlobal.Range("C2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Give the new column a sensible width (~9.42 characters).
$wsGlobal.Columns.Item(4).ColumnWidth = 8.6

# Remember the new selection on the Global sheet without changing which
# sheet/tab is actually active in the workbook.
$wsGlobal.Activate()
$wsGlobal.Range("D2").Select()
$wsActive.Activate()
